# Corrigindo para apenas deputados federais de 2018
# Insert a new "Foto" column before the existing "SQ_CANDIDATO" column (D),
# shifting the old D:I columns to E:J, and fill the new column with the
# photo path derived from the (now shifted) SQ_CANDIDATO value in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; this shifts old D..I to E..J and
# automatically grows the sheet dimension/used range (A1:I74 -> A1:J74).
$ws.Columns("D").Insert()

# Header for the newly inserted column.
$ws.Range("D1").Value = "Foto"

# Fill in the photo reference for every data row, based on SQ_CANDIDATO
# (now located in column E after the insert).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $sq = $ws.Cells.Item($r, 5).Value2
    $padded = $sq.ToString().PadLeft(10, '0')
    $ws.Cells.Item($r, 4).Value = " foto_cand2014_div/FBR28" + $padded
}
